# Added function declaration for toggle and radio buttons.
# Two new auto-generated translation text entries ("SingleUseId52" and
# "SingleUseId53") were inserted into the TouchGFX text-id sequence,
# which cascades a renumbering of the ids used by rows 40-53 of the
# "Translation" sheet (column B), along with their associated
# alignment/typography and resource-text values in columns C/D/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B40").Value = "SingleUseId39"
$ws.Range("D40").Value = "Left"
$ws.Range("F40").Value = "Meas Rate"
$ws.Range("B41").Value = "SingleUseId40"
$ws.Range("D41").Value = "Center"
$ws.Range("F41").Value = "<value> ms"
$ws.Range("B42").Value = "SingleUseId42"
$ws.Range("F42").Value = "<value>"
$ws.Range("B43").Value = "SingleUseId43"
$ws.Range("D43").Value = "Left"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "0"
$ws.Range("F43").Style = "Normal"
$ws.Range("B44").Value = "SingleUseId44"
$ws.Range("D44").Value = "Left"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "0"
$ws.Range("F44").Style = "Normal"
$ws.Range("B45").Value = "SingleUseId45"
$ws.Range("F45").Value = "Stamps Number"
$ws.Range("B46").Value = "SingleUseId46"
$ws.Range("C46").Value = "Large"
$ws.Range("F46").Value = "X"
$ws.Range("B47").Value = "SingleUseId47"
$ws.Range("F47").Value = "Repeat"
$ws.Range("B48").Value = "SingleUseId48"
$ws.Range("C48").Value = "Default"
$ws.Range("D48").Value = "Center"
$ws.Range("F48").Value = "<value>"
$ws.Range("B49").Value = "SingleUseId49"
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value = "0"
$ws.Range("F49").Style = "Normal"
$ws.Range("B50").Value = "SingleUseId50"
$ws.Range("F50").Value = "<value> s"
$ws.Range("B51").Value = "SingleUseId51"
$ws.Range("B52").Value = "SingleUseId52"
$ws.Range("D52").Value = "Left"
$ws.Range("F52").Value = "Single"
$ws.Range("B53").Value = "SingleUseId53"
$ws.Range("F53").Value = "Continuous"

Write-Output "Applied translation sheet updates for rows 40-53"
